# Weekly update: insert two new price rows (Primera / Segunda) for Cebollín at
# Vega Monumental Concepción, pushing the existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 31 (first data row after the
# header and the first 29 existing records that stay untouched).
$ws.Rows("31:32").Insert()

# --- New row 31 : "Primera" quality record, dated 2023-02-08 (serial 44965) ---
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 44965
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112037
$ws.Range("G31").Value = "Cebollín"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 700
$ws.Range("L31").Value = 800
$ws.Range("M31").Value = 750
$ws.Range("N31").Value = "`$/paquete 6 unidades"
$ws.Range("O31").Value = "Región de Ñuble"
$ws.Range("P31").Value = 125
$ws.Range("Q31").Value = 6
$ws.Range("R31").Value = "Hortaliza"

# --- New row 32 : "Segunda" quality record, same date ---
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44965
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112037
$ws.Range("G32").Value = "Cebollín"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = 600
$ws.Range("N32").Value = "`$/paquete 6 unidades"
$ws.Range("O32").Value = "Región de Ñuble"
$ws.Range("P32").Value = 100
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = "Hortaliza"

# Apply the same date-number format used by the other "Fecha" cells (s="2" in the
# original workbook) so the new D31/D32 cells render consistently.
$ws.Range("D31:D32").NumberFormat = $ws.Range("D33").NumberFormat
